# [Kadastro App] Yeni kayit eklendi: 3007
#
# Adds one new "Kayit No 3007" record row to the bottom of the two
# worksheets that list it: the "Kayitlar" master table and the
# "Erdemli" per-office table (both currently end at row 66).

$wb = $excel.ActiveWorkbook

# Kayıt No, Tarih, Birim, Parsel Sayısı, İş, Personeller
$kayitNo    = "3007"
$tarih      = "2025-09-11"
$birim      = "Erdemli"
$parselSayi = "1"
$is         = "CİNS DEĞ."
$personel   = "CEMAL TİMUROĞLU (K.Teknisyeni), ÖZKAN AKBAŞ (Mühendis)"

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the next empty row right below the existing data (column A).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    $rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 6))
    # Keep the numeric-looking fields (Kayıt No / Tarih / Parsel Sayısı) as
    # text, matching every other row in this sheet.
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($newRow, 1).Value = $kayitNo
    $ws.Cells.Item($newRow, 2).Value = $tarih
    $ws.Cells.Item($newRow, 3).Value = $birim
    $ws.Cells.Item($newRow, 4).Value = $parselSayi
    $ws.Cells.Item($newRow, 5).Value = $is
    $ws.Cells.Item($newRow, 6).Value = $personel
}
